$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9080.1755268314
$ws.Range("C2").Value = 8649.30968855622
$ws.Range("E2").Value = 4881.23206108236
$ws.Range("F2").Value = 11.6059062349411

# Row 3
$ws.Range("B3").Value = 8863.21007302798
$ws.Range("C3").Value = 8393.47708863
$ws.Range("E3").Value = 4761.15119828872
$ws.Range("F3").Value = 172.94284528828

# Row 4
$ws.Range("B4").Value = 8828.20494188129
$ws.Range("C4").Value = 8231.60835028267
$ws.Range("E4").Value = 4926.43029401483
$ws.Range("F4").Value = 173.084943512396

# Row 5
$ws.Range("B5").Value = 8761.40850209054
$ws.Range("C5").Value = 7341.34743769629
$ws.Range("E5").Value = 4871.18940750456
$ws.Range("F5").Value = 133.689035216702

# Row 6
$ws.Range("B6").Value = 3074.34101973407
$ws.Range("C6").Value = 4864.8010757703
$ws.Range("E6").Value = 4378.08489325
$ws.Range("F6").Value = 9.95358204251238

# Row 7
$ws.Range("B7").Value = 2850.36187997388
$ws.Range("C7").Value = 4498.68111762163
$ws.Range("E7").Value = 4266.66384808243
$ws.Range("F7").Value = -9.94395976233075

# Row 9
$ws.Range("C9").Value = 8754.81719730593
$ws.Range("F9").Value = 232.508871742085

# Row 10
$ws.Range("C10").Value = 9295.32663079935
$ws.Range("F10").Value = 255.030098137645

# Row 11
$ws.Range("C11").Value = 9954.81529540331
$ws.Range("F11").Value = 282.508792496143

# Row 12
$ws.Range("C12").Value = 9484.94736584209
$ws.Range("F12").Value = 262.930962097759

# Row 13
$ws.Range("C13").Value = 6984.97646269771
$ws.Range("F13").Value = 144.105546836205

# Row 14
$ws.Range("C14").Value = 6818.03462064937
$ws.Range("F14").Value = 136.816157608522

# Row 15
$ws.Range("C15").Value = 10694.2187751014
$ws.Range("F15").Value = 337.46884561394
